# Update numeric values in the second table ("Regression") of the
# document. Table 1 (Classification) is left untouched; only Table 2
# (Regression) has cell values changed, per the commit's updated figures.

$d = $word.ActiveDocument
$t = $d.Tables.Item(2)

# Each entry: row, column (1-based, as in Word's Table.Cell(row, col)),
# expected old text, new text.
$changes = @(
    @(2, 2, "0.199", "0.187"),
    @(2, 3, "0.026", "0.024"),
    @(2, 4, "0.416", "0.417"),

    @(3, 2, "0.091", "0.097"),
    @(3, 3, "0.022", "0.024"),
    @(3, 4, "0.317", "0.319"),
    @(3, 5, "0.026", "0.027"),

    @(4, 2, "0.015", "0.005"),
    @(4, 3, "0.007", "0.006"),
    @(4, 4, "0.337", "0.338"),
    @(4, 5, "0.023", "0.024"),

    @(5, 2, "0.119", "0.126"),
    @(5, 3, "0.014", "0.019"),
    @(5, 4, "0.124", "0.125"),
    @(5, 6, "0.833", "0.973"),

    @(6, 2, "0.132", "0.135"),
    @(6, 3, "0.026", "0.022"),
    @(6, 4, "0.324", "0.325"),

    @(7, 2, "0.157", "0.155"),
    @(7, 4, "0.383", "0.381"),
    @(7, 5, "0.022", "0.021"),

    @(8, 2, "0.172", "0.145"),
    @(8, 3, "0.030", "0.022"),
    @(8, 4, "0.329", "0.332"),
    @(8, 5, "0.020", "0.019"),
    @(8, 6, "0.001", "0.000")
)

foreach ($change in $changes) {
    $row = $change[0]
    $col = $change[1]
    $old = $change[2]
    $new = $change[3]

    $cell = $t.Cell($row, $col)
    $current = $cell.Range.Text.TrimEnd([char]13, [char]7)

    if ($current -ne $old) {
        Write-Host "Mismatch at row $row col $col : expected '$old' found '$current'"
    }

    $cell.Range.Text = $new
}
